# Update sensitivity results (commit: "update sensitivity results final 0606")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "I5" = 2259

    "D6" = 109
    "E6" = 93
    "F6" = 0
    "G6" = 108
    "H6" = 93
    "I6" = 0
    "K6" = 45
    "L6" = 105
    "M6" = 105
    "N6" = 108
    "O6" = 0

    "E7" = 0
    "G7" = 172
    "H7" = 0
    "I7" = 0
    "K7" = 18
    "M7" = 172
    "N7" = 172
    "O7" = 0

    "D8" = 1089
    "E8" = 1304
    "F8" = 0
    "G8" = 1083
    "H8" = 0
    "I8" = 0
    "K8" = 1092
    "L8" = 162
    "M8" = 1066
    "N8" = 0
    "O8" = 0

    "E9" = 0
    "H9" = 1303
    "I9" = 738
    "K9" = 0
    "N9" = 1066
    "O9" = 1067

    "F10" = 56
    "I10" = 34
    "L10" = 56
    "O10" = 75

    "G11" = 157
    "H11" = 157
    "I11" = 157
    "M11" = 157
    "N11" = 157
    "O11" = 157
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
